# Apply the "data error" test-case additions for Shipping Details / Shipping
# Address / Duplicate PO stories to the Input_Data worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Three new header columns are appended after column AK:
#    AL - X4DDataErrorOrderConfirmationId
#    AM - ModifyShippingAddressDataErrorOrderID
#    AN - X4CDuplicatePODataErrorOrderID
#    Copy the formatting of the last existing header cell (AK1) onto
#    the new header cells first, then fill in their text.
# ------------------------------------------------------------------
$ws.Range("AK1").Copy()
$ws.Range("AL1:AN1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AL1").Value = "X4DDataErrorOrderConfirmationId"
$ws.Range("AM1").Value = "ModifyShippingAddressDataErrorOrderID"

# ------------------------------------------------------------------
# 2. Row 5 (the "A5" label) is renamed from the old data-error bucket
#    name to the new "order exception" bucket name.
# ------------------------------------------------------------------
$ws.Range("A5").Value = "order_exception_orders"

$ws.Range("AN1").Value = "X4CDuplicatePODataErrorOrderID"

# ------------------------------------------------------------------
# 3. Column width adjustments.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.333333333333336
$ws.Columns.Item(38).ColumnWidth = 32.83333333333333
$ws.Columns.Item(39).ColumnWidth = 43.5
$ws.Columns.Item(40).ColumnWidth = 31.5

# ------------------------------------------------------------------
# 4. Update the view so the new columns are visible / selected,
#    matching the author's final cursor position.
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 38
$win.ScrollRow = 1
$ws.Range("AP5").Select()
